$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.106.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.358.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +10.11%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.98"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +10.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "621.86"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.61%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +10.76%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.353.11"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +10.03%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "97.833.76"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.71"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.80%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.978.67"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +10.01%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000245"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.355.53"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +9.79%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.98"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "483.95"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +11.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.82"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000205"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +8.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.64"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.28"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +10.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.253"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.185"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.88%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -9.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.25"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.20"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.98%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "517.86"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +12.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.150"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.03%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.60"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.96%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.78%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.774"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +16.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.58"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.64%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +8.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.50"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.08%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.47%  "
